$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: write one "history" row consisting of
#   A: run date (text, e.g. 2023-04-19)
#   B: run time (datetime serial, formatted as YYYY-MM-DD HH:MM:SS)
#   C: sprint / run name (text)
#   D: total cases (number)
#   E: pass cases (number)
#   F: fail cases (number)
#   G: time taken (number)
# ---------------------------------------------------------------
function Set-HistoryRow($ws, $row, $dateText, $timeSerial, $nameText, $total, $pass, $fail, $taken) {
    # Force column A & C to remain plain text (leading apostrophe keeps
    # Excel from auto-converting the date-looking string into a date).
    $ws.Cells.Item($row, 1).Formula = "'" + $dateText
    $ws.Cells.Item($row, 3).Formula = "'" + $nameText

    # Column B keeps the existing "YYYY-MM-DD HH:MM:SS" custom format
    # used throughout the sheet for the run timestamp.
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Formula = "$timeSerial"

    $ws.Cells.Item($row, 4).Formula = "$total"
    $ws.Cells.Item($row, 5).Formula = "$pass"
    $ws.Cells.Item($row, 6).Formula = "$fail"
    $ws.Cells.Item($row, 7).Formula = "$taken"
}

# ---------------------------------------------------------------
# AMSIN sheet - new tenant rows for second cycle + final run
# ---------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Existing row 71 picks up the refreshed run timestamp captured when the
# two additional rows below it were produced.
$wsAmsin.Cells.Item(71, 2).Formula = "45034.67370373842"

Set-HistoryRow $wsAmsin 72 "2023-04-19" "45035.69311298611" "176scndcyc" 124 124 0 1.83
Set-HistoryRow $wsAmsin 73 "2023-04-20" "45036.40967388889" "176fnlruntest" 124 124 0 1.68

# ---------------------------------------------------------------
# BETA sheet - new tenant row
# ---------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
Set-HistoryRow $wsBeta 35 "2023-04-20" "45036.51260524306" "176beta" 124 122 2 1.61

# ---------------------------------------------------------------
# AMS sheet - new tenant row
# ---------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")
Set-HistoryRow $wsAms 39 "2023-05-08" "45054.53683669174" "176htfxtrl" 124 123 1 1.72
